$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row right above the "001882235" account row (currently row 2)
#    and populate it with the new account data (EVANGELINA).
$target = $ws.Columns.Item(1).Find("001882235")
$insertRow = $target.Row
$ws.Rows.Item($insertRow).Insert()

$ws.Cells.Item($insertRow, 1).Value = "'005646524"
$ws.Cells.Item($insertRow, 2).Value = "EVANGELINA"
$ws.Cells.Item($insertRow, 3).Value = 1000000

# 2. Remove the row for account 005024046 (ALEXANDRE)
$row1 = $ws.Columns.Item(1).Find("005024046")
$ws.Rows.Item($row1.Row).Delete()

# 3. Remove the row for account 004752461 (SERGIO)
$row2 = $ws.Columns.Item(1).Find("004752461")
$ws.Rows.Item($row2.Row).Delete()
